$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "Trinkbecher getroffen" column header (AA1), matching header style ---
$ws.Range("Z1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)
$ws.Range("AA1").Value = "Trinkbecher getroffen"

# --- Apply the existing date-column (I) number format to the new rows before filling values ---
$ws.Range("I244").Copy()
$ws.Range("I245:I269").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new match results (rows 245-269) ---
$ws.Cells.Item(245, 1).Value = "Philipp"
$ws.Cells.Item(245, 2).Value = "André"
$ws.Cells.Item(245, 4).Value = "Marie"
$ws.Cells.Item(245, 5).Value = "Emilio"
$ws.Cells.Item(245, 7).Value = "Garruk Ultras"
$ws.Cells.Item(245, 8).Value = "Die Prenzlauer Crew"
$ws.Cells.Item(245, 9).Value = 45479
$ws.Cells.Item(245, 10).Value = 1
$ws.Cells.Item(245, 11).Value = 0
$ws.Cells.Item(245, 12).Value = 1
$ws.Cells.Item(245, 14).Value = 1
$ws.Cells.Item(245, 15).Value = 1
$ws.Cells.Item(245, 16).Value = 6
$ws.Cells.Item(245, 17).Value = 3
$ws.Cells.Item(245, 19).Value = 3
$ws.Cells.Item(245, 20).Value = 5
$ws.Cells.Item(246, 1).Value = "Chris"
$ws.Cells.Item(246, 2).Value = "Valdemar"
$ws.Cells.Item(246, 4).Value = "Phia"
$ws.Cells.Item(246, 5).Value = "Daisy"
$ws.Cells.Item(246, 7).Value = "Pferdewelle Stampf Stampf"
$ws.Cells.Item(246, 8).Value = "Fleißigen Bienen Bzz Bzz"
$ws.Cells.Item(246, 9).Value = 45479
$ws.Cells.Item(246, 10).Value = 1
$ws.Cells.Item(246, 11).Value = 0
$ws.Cells.Item(246, 12).Value = 2
$ws.Cells.Item(246, 13).Value = "Chris"
$ws.Cells.Item(246, 14).Value = 1
$ws.Cells.Item(246, 16).Value = 5
$ws.Cells.Item(246, 17).Value = 1
$ws.Cells.Item(246, 19).Value = 2
$ws.Cells.Item(246, 20).Value = 2
$ws.Cells.Item(247, 1).Value = "Phia"
$ws.Cells.Item(247, 2).Value = "Daisy"
$ws.Cells.Item(247, 4).Value = "Philipp"
$ws.Cells.Item(247, 5).Value = "André"
$ws.Cells.Item(247, 7).Value = "Fleißigen Bienen Bzz Bzz"
$ws.Cells.Item(247, 8).Value = "Garruk Ultras"
$ws.Cells.Item(247, 9).Value = 45479
$ws.Cells.Item(247, 10).Value = 0
$ws.Cells.Item(247, 11).Value = 1
$ws.Cells.Item(247, 12).Value = 2
$ws.Cells.Item(247, 14).Value = 1
$ws.Cells.Item(247, 16).Value = 0
$ws.Cells.Item(247, 17).Value = 4
$ws.Cells.Item(247, 19).Value = 4
$ws.Cells.Item(247, 20).Value = 2
$ws.Cells.Item(248, 1).Value = "Marie"
$ws.Cells.Item(248, 2).Value = "Emilio"
$ws.Cells.Item(248, 4).Value = "Chris"
$ws.Cells.Item(248, 5).Value = "Valdemar"
$ws.Cells.Item(248, 7).Value = "Die Prenzlauer Crew"
$ws.Cells.Item(248, 8).Value = "Pferdewelle Stampf Stampf"
$ws.Cells.Item(248, 9).Value = 45479
$ws.Cells.Item(248, 10).Value = 1
$ws.Cells.Item(248, 11).Value = 0
$ws.Cells.Item(248, 12).Value = 1
$ws.Cells.Item(248, 14).Value = 1
$ws.Cells.Item(248, 16).Value = 3
$ws.Cells.Item(248, 17).Value = 3
$ws.Cells.Item(248, 19).Value = 3
$ws.Cells.Item(248, 20).Value = 2
$ws.Cells.Item(249, 1).Value = "Phia"
$ws.Cells.Item(249, 2).Value = "Daisy"
$ws.Cells.Item(249, 4).Value = "Marie"
$ws.Cells.Item(249, 5).Value = "Emilio"
$ws.Cells.Item(249, 7).Value = "Fleißigen Bienen Bzz Bzz"
$ws.Cells.Item(249, 8).Value = "Die Prenzlauer Crew"
$ws.Cells.Item(249, 9).Value = 45479
$ws.Cells.Item(249, 10).Value = 1
$ws.Cells.Item(249, 11).Value = 0
$ws.Cells.Item(249, 12).Value = 1
$ws.Cells.Item(249, 14).Value = 1
$ws.Cells.Item(249, 16).Value = 2
$ws.Cells.Item(249, 17).Value = 4
$ws.Cells.Item(249, 19).Value = 1
$ws.Cells.Item(249, 20).Value = 4
$ws.Cells.Item(250, 1).Value = "Chris"
$ws.Cells.Item(250, 2).Value = "Valdemar"
$ws.Cells.Item(250, 4).Value = "Philipp"
$ws.Cells.Item(250, 5).Value = "André"
$ws.Cells.Item(250, 7).Value = "Pferdewelle Stampf Stampf"
$ws.Cells.Item(250, 8).Value = "Garruk Ultras"
$ws.Cells.Item(250, 9).Value = 45479
$ws.Cells.Item(250, 10).Value = 0
$ws.Cells.Item(250, 11).Value = 1
$ws.Cells.Item(250, 12).Value = 2
$ws.Cells.Item(250, 14).Value = 1
$ws.Cells.Item(250, 16).Value = 1
$ws.Cells.Item(250, 17).Value = 3
$ws.Cells.Item(250, 19).Value = 4
$ws.Cells.Item(250, 20).Value = 2
$ws.Cells.Item(251, 1).Value = "Chris"
$ws.Cells.Item(251, 2).Value = "Phia"
$ws.Cells.Item(251, 4).Value = "Emilio"
$ws.Cells.Item(251, 5).Value = "Philipp"
$ws.Cells.Item(251, 6).Value = "Daisy"
$ws.Cells.Item(251, 7).Value = "DD - Dirk & Dora"
$ws.Cells.Item(251, 8).Value = "Alles Andere als Arbeit"
$ws.Cells.Item(251, 9).Value = 45479
$ws.Cells.Item(251, 10).Value = 1
$ws.Cells.Item(251, 11).Value = 0
$ws.Cells.Item(251, 12).Value = 4
$ws.Cells.Item(251, 14).Value = 1
$ws.Cells.Item(251, 16).Value = 2
$ws.Cells.Item(251, 17).Value = 4
$ws.Cells.Item(251, 19).Value = 1
$ws.Cells.Item(251, 20).Value = 1
$ws.Cells.Item(251, 21).Value = 0
$ws.Cells.Item(252, 1).Value = "Emilio"
$ws.Cells.Item(252, 2).Value = "Philipp"
$ws.Cells.Item(252, 3).Value = "Daisy"
$ws.Cells.Item(252, 4).Value = "Valdemar"
$ws.Cells.Item(252, 5).Value = "André"
$ws.Cells.Item(252, 7).Value = "Alles Andere als Arbeit"
$ws.Cells.Item(252, 8).Value = "Zuckerlager voll"
$ws.Cells.Item(252, 9).Value = 45479
$ws.Cells.Item(252, 10).Value = 0
$ws.Cells.Item(252, 11).Value = 1
$ws.Cells.Item(252, 12).Value = 2
$ws.Cells.Item(252, 13).Value = "Emilio, André"
$ws.Cells.Item(252, 14).Value = 1
$ws.Cells.Item(252, 15).Value = 1
$ws.Cells.Item(252, 16).Value = 2
$ws.Cells.Item(252, 17).Value = 3
$ws.Cells.Item(252, 18).Value = 2
$ws.Cells.Item(252, 19).Value = 5
$ws.Cells.Item(252, 20).Value = 4
$ws.Cells.Item(253, 1).Value = "Valdemar"
$ws.Cells.Item(253, 2).Value = "André"
$ws.Cells.Item(253, 4).Value = "Chris"
$ws.Cells.Item(253, 5).Value = "Phia"
$ws.Cells.Item(253, 7).Value = "Zuckerlager voll"
$ws.Cells.Item(253, 8).Value = "DD - Dirk & Dora"
$ws.Cells.Item(253, 9).Value = 45479
$ws.Cells.Item(253, 10).Value = 1
$ws.Cells.Item(253, 11).Value = 0
$ws.Cells.Item(253, 12).Value = 2
$ws.Cells.Item(253, 13).Value = "André, André"
$ws.Cells.Item(253, 14).Value = 1
$ws.Cells.Item(253, 16).Value = 3
$ws.Cells.Item(253, 17).Value = 3
$ws.Cells.Item(253, 19).Value = 2
$ws.Cells.Item(253, 20).Value = 2
$ws.Cells.Item(254, 1).Value = "Phia"
$ws.Cells.Item(254, 2).Value = "André"
$ws.Cells.Item(254, 4).Value = "Valdemar"
$ws.Cells.Item(254, 5).Value = "Emilio"
$ws.Cells.Item(254, 7).Value = "Team Heul doch"
$ws.Cells.Item(254, 8).Value = "Die romantischen Matrosen"
$ws.Cells.Item(254, 9).Value = 45479
$ws.Cells.Item(254, 10).Value = 1
$ws.Cells.Item(254, 11).Value = 0
$ws.Cells.Item(254, 12).Value = 1
$ws.Cells.Item(254, 14).Value = 1
$ws.Cells.Item(254, 16).Value = 2
$ws.Cells.Item(254, 17).Value = 4
$ws.Cells.Item(254, 19).Value = 4
$ws.Cells.Item(254, 20).Value = 1
$ws.Cells.Item(255, 1).Value = "Chris"
$ws.Cells.Item(255, 2).Value = "Phia"
$ws.Cells.Item(255, 4).Value = "Valdemar"
$ws.Cells.Item(255, 5).Value = "Emilio"
$ws.Cells.Item(255, 7).Value = "ZaZa Grill"
$ws.Cells.Item(255, 8).Value = "I got bit by a WIDDER"
$ws.Cells.Item(255, 9).Value = 45500
$ws.Cells.Item(255, 10).Value = 0
$ws.Cells.Item(255, 11).Value = 1
$ws.Cells.Item(255, 12).Value = 1
$ws.Cells.Item(255, 14).Value = 3
$ws.Cells.Item(255, 16).Value = 3
$ws.Cells.Item(255, 17).Value = 2
$ws.Cells.Item(255, 19).Value = 2
$ws.Cells.Item(255, 20).Value = 4
$ws.Cells.Item(256, 1).Value = "André"
$ws.Cells.Item(256, 2).Value = "Marian"
$ws.Cells.Item(256, 4).Value = "Chris"
$ws.Cells.Item(256, 5).Value = "Phia"
$ws.Cells.Item(256, 7).Value = "Bowle Batallion"
$ws.Cells.Item(256, 8).Value = "ZaZa Grill"
$ws.Cells.Item(256, 9).Value = 45500
$ws.Cells.Item(256, 10).Value = 0
$ws.Cells.Item(256, 11).Value = 1
$ws.Cells.Item(256, 12).Value = 2
$ws.Cells.Item(256, 14).Value = 3
$ws.Cells.Item(256, 16).Value = 1
$ws.Cells.Item(256, 17).Value = 3
$ws.Cells.Item(256, 19).Value = 3
$ws.Cells.Item(256, 20).Value = 3
$ws.Cells.Item(257, 1).Value = "Valdemar"
$ws.Cells.Item(257, 2).Value = "Emilio"
$ws.Cells.Item(257, 4).Value = "André"
$ws.Cells.Item(257, 5).Value = "Marian"
$ws.Cells.Item(257, 7).Value = "I got bit by a WIDDER"
$ws.Cells.Item(257, 8).Value = "Bowle Batallion"
$ws.Cells.Item(257, 9).Value = 45500
$ws.Cells.Item(257, 10).Value = 1
$ws.Cells.Item(257, 11).Value = 0
$ws.Cells.Item(257, 12).Value = 1
$ws.Cells.Item(257, 14).Value = 3
$ws.Cells.Item(257, 16).Value = 2
$ws.Cells.Item(257, 17).Value = 4
$ws.Cells.Item(257, 19).Value = 2
$ws.Cells.Item(257, 20).Value = 3
$ws.Cells.Item(258, 1).Value = "Emilio"
$ws.Cells.Item(258, 2).Value = "Phia"
$ws.Cells.Item(258, 4).Value = "André"
$ws.Cells.Item(258, 5).Value = "Leonie"
$ws.Cells.Item(258, 7).Value = "K-Hole"
$ws.Cells.Item(258, 8).Value = "Schnelle Bälle"
$ws.Cells.Item(258, 9).Value = 45500
$ws.Cells.Item(258, 10).Value = 1
$ws.Cells.Item(258, 11).Value = 0
$ws.Cells.Item(258, 12).Value = 1
$ws.Cells.Item(258, 14).Value = 3
$ws.Cells.Item(258, 16).Value = 3
$ws.Cells.Item(258, 17).Value = 3
$ws.Cells.Item(258, 19).Value = 0
$ws.Cells.Item(258, 20).Value = 5
$ws.Cells.Item(259, 1).Value = "Chris"
$ws.Cells.Item(259, 2).Value = "Marian"
$ws.Cells.Item(259, 4).Value = "Merlin"
$ws.Cells.Item(259, 5).Value = "Valdemar"
$ws.Cells.Item(259, 7).Value = "Best LoL-Players in the room"
$ws.Cells.Item(259, 8).Value = "Two Bikey Boys Go Vroom"
$ws.Cells.Item(259, 9).Value = 45500
$ws.Cells.Item(259, 10).Value = 0
$ws.Cells.Item(259, 11).Value = 1
$ws.Cells.Item(259, 12).Value = 2
$ws.Cells.Item(259, 14).Value = 3
$ws.Cells.Item(259, 16).Value = 3
$ws.Cells.Item(259, 17).Value = 1
$ws.Cells.Item(259, 19).Value = 3
$ws.Cells.Item(259, 20).Value = 3
$ws.Cells.Item(260, 1).Value = "Merlin"
$ws.Cells.Item(260, 2).Value = "Valdemar"
$ws.Cells.Item(260, 4).Value = "Emilio"
$ws.Cells.Item(260, 5).Value = "Phia"
$ws.Cells.Item(260, 7).Value = "Two Bikey Boys Go Vroom"
$ws.Cells.Item(260, 8).Value = "K-Hole"
$ws.Cells.Item(260, 9).Value = 45500
$ws.Cells.Item(260, 10).Value = 0
$ws.Cells.Item(260, 11).Value = 1
$ws.Cells.Item(260, 12).Value = 1
$ws.Cells.Item(260, 14).Value = 3
$ws.Cells.Item(260, 16).Value = 2
$ws.Cells.Item(260, 17).Value = 3
$ws.Cells.Item(260, 19).Value = 3
$ws.Cells.Item(260, 20).Value = 3
$ws.Cells.Item(261, 1).Value = "André"
$ws.Cells.Item(261, 2).Value = "Leonie"
$ws.Cells.Item(261, 4).Value = "Chris"
$ws.Cells.Item(261, 5).Value = "Marian"
$ws.Cells.Item(261, 7).Value = "Schnelle Bälle"
$ws.Cells.Item(261, 8).Value = "Best LoL-Players in the room"
$ws.Cells.Item(261, 9).Value = 45500
$ws.Cells.Item(261, 10).Value = 0
$ws.Cells.Item(261, 11).Value = 1
$ws.Cells.Item(261, 12).Value = 1
$ws.Cells.Item(261, 14).Value = 3
$ws.Cells.Item(261, 16).Value = 2
$ws.Cells.Item(261, 17).Value = 3
$ws.Cells.Item(261, 19).Value = 5
$ws.Cells.Item(261, 20).Value = 1
$ws.Cells.Item(262, 1).Value = "André"
$ws.Cells.Item(262, 2).Value = "Leonie"
$ws.Cells.Item(262, 4).Value = "Merlin"
$ws.Cells.Item(262, 5).Value = "Valdemar"
$ws.Cells.Item(262, 7).Value = "Schnelle Bälle"
$ws.Cells.Item(262, 8).Value = "Two Bikey Boys Go Vroom"
$ws.Cells.Item(262, 9).Value = 45500
$ws.Cells.Item(262, 10).Value = 1
$ws.Cells.Item(262, 11).Value = 0
$ws.Cells.Item(262, 12).Value = 2
$ws.Cells.Item(262, 14).Value = 3
$ws.Cells.Item(262, 16).Value = 2
$ws.Cells.Item(262, 17).Value = 4
$ws.Cells.Item(262, 19).Value = 3
$ws.Cells.Item(262, 20).Value = 1
$ws.Cells.Item(263, 1).Value = "Emilio"
$ws.Cells.Item(263, 2).Value = "Phia"
$ws.Cells.Item(263, 4).Value = "Chris"
$ws.Cells.Item(263, 5).Value = "Marian"
$ws.Cells.Item(263, 7).Value = "K-Hole"
$ws.Cells.Item(263, 8).Value = "Best LoL-Players in the room"
$ws.Cells.Item(263, 9).Value = 45500
$ws.Cells.Item(263, 10).Value = 1
$ws.Cells.Item(263, 11).Value = 0
$ws.Cells.Item(263, 12).Value = 1
$ws.Cells.Item(263, 14).Value = 3
$ws.Cells.Item(263, 16).Value = 3
$ws.Cells.Item(263, 17).Value = 3
$ws.Cells.Item(263, 19).Value = 1
$ws.Cells.Item(263, 20).Value = 4
$ws.Cells.Item(264, 1).Value = "André"
$ws.Cells.Item(264, 2).Value = "Emilio"
$ws.Cells.Item(264, 4).Value = "Merlin"
$ws.Cells.Item(264, 5).Value = "Leonie"
$ws.Cells.Item(264, 7).Value = "Pokerogue Prodigies"
$ws.Cells.Item(264, 8).Value = "So ein großer Yarak"
$ws.Cells.Item(264, 9).Value = 45500
$ws.Cells.Item(264, 10).Value = 0
$ws.Cells.Item(264, 11).Value = 0
$ws.Cells.Item(264, 12).Value = 1
$ws.Cells.Item(264, 14).Value = 3
$ws.Cells.Item(264, 16).Value = 3
$ws.Cells.Item(264, 17).Value = 2
$ws.Cells.Item(264, 19).Value = 4
$ws.Cells.Item(264, 20).Value = 1
$ws.Cells.Item(265, 1).Value = "Marian"
$ws.Cells.Item(265, 2).Value = "Phia"
$ws.Cells.Item(265, 4).Value = "Valdemar"
$ws.Cells.Item(265, 5).Value = "Chris"
$ws.Cells.Item(265, 7).Value = "Redemption Arc"
$ws.Cells.Item(265, 8).Value = "Mehmet Arms"
$ws.Cells.Item(265, 9).Value = 45500
$ws.Cells.Item(265, 10).Value = 0
$ws.Cells.Item(265, 11).Value = 1
$ws.Cells.Item(265, 12).Value = 1
$ws.Cells.Item(265, 14).Value = 3
$ws.Cells.Item(265, 16).Value = 2
$ws.Cells.Item(265, 17).Value = 3
$ws.Cells.Item(265, 19).Value = 4
$ws.Cells.Item(265, 20).Value = 2
$ws.Cells.Item(266, 1).Value = "Marian"
$ws.Cells.Item(266, 2).Value = "Phia"
$ws.Cells.Item(266, 4).Value = "Merlin"
$ws.Cells.Item(266, 5).Value = "Leonie"
$ws.Cells.Item(266, 7).Value = "Redemption Arc"
$ws.Cells.Item(266, 8).Value = "So ein großer Yarak"
$ws.Cells.Item(266, 9).Value = 45500
$ws.Cells.Item(266, 10).Value = 1
$ws.Cells.Item(266, 11).Value = 0
$ws.Cells.Item(266, 12).Value = 1
$ws.Cells.Item(266, 14).Value = 3
$ws.Cells.Item(266, 16).Value = 2
$ws.Cells.Item(266, 17).Value = 4
$ws.Cells.Item(266, 19).Value = 3
$ws.Cells.Item(266, 20).Value = 2
$ws.Cells.Item(267, 1).Value = "Valdemar"
$ws.Cells.Item(267, 2).Value = "Chris"
$ws.Cells.Item(267, 4).Value = "André"
$ws.Cells.Item(267, 5).Value = "Emilio"
$ws.Cells.Item(267, 7).Value = "Mehmet Arms"
$ws.Cells.Item(267, 8).Value = "Pokerogue Prodigies"
$ws.Cells.Item(267, 9).Value = 45500
$ws.Cells.Item(267, 10).Value = 0
$ws.Cells.Item(267, 11).Value = 1
$ws.Cells.Item(267, 12).Value = 2
$ws.Cells.Item(267, 14).Value = 3
$ws.Cells.Item(267, 16).Value = 1
$ws.Cells.Item(267, 17).Value = 3
$ws.Cells.Item(267, 19).Value = 2
$ws.Cells.Item(267, 20).Value = 4
$ws.Cells.Item(268, 1).Value = "André"
$ws.Cells.Item(268, 2).Value = "Emilio"
$ws.Cells.Item(268, 4).Value = "Marian"
$ws.Cells.Item(268, 5).Value = "Phia"
$ws.Cells.Item(268, 7).Value = "Pokerogue Prodigies"
$ws.Cells.Item(268, 8).Value = "Redemption Arc"
$ws.Cells.Item(268, 9).Value = 45500
$ws.Cells.Item(268, 10).Value = 0
$ws.Cells.Item(268, 11).Value = 1
$ws.Cells.Item(268, 12).Value = 4
$ws.Cells.Item(268, 14).Value = 3
$ws.Cells.Item(268, 16).Value = 3
$ws.Cells.Item(268, 17).Value = 3
$ws.Cells.Item(268, 19).Value = 1
$ws.Cells.Item(268, 20).Value = 1
$ws.Cells.Item(269, 1).Value = "Merlin"
$ws.Cells.Item(269, 2).Value = "Leonie"
$ws.Cells.Item(269, 4).Value = "Valdemar"
$ws.Cells.Item(269, 5).Value = "Chris"
$ws.Cells.Item(269, 7).Value = "So ein großer Yarak"
$ws.Cells.Item(269, 8).Value = "Mehmet Arms"
$ws.Cells.Item(269, 9).Value = 45500
$ws.Cells.Item(269, 10).Value = 0
$ws.Cells.Item(269, 11).Value = 1
$ws.Cells.Item(269, 12).Value = 3
$ws.Cells.Item(269, 14).Value = 3
$ws.Cells.Item(269, 16).Value = 3
$ws.Cells.Item(269, 19).Value = 1
$ws.Cells.Item(269, 20).Value = 3
$ws.Cells.Item(269, 27).Value = "Chris -> Merlin"


# --- Restore selection/active cell to match the final editing state ---
$ws.Range("T269").Select() | Out-Null
